$wb = $excel.ActiveWorkbook
$wsTablet = $wb.Worksheets.Item("Tablet")
$co = $wsTablet.ChartObjects().Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()
$s1 = $sc.Item(1)
$members = [System.__ComObject].GetType()
Write-Host ($s1.GetType().FullName)
